{"js": "// Rettet en smule i Resum\u00e9.docx\n// Applies the wording tweaks to the two \"Pristjek220\" summary paragraphs and\n// relocates the `_GoBack` bookmark from the end of the \"Gruppen har arbejdet...\"\n// paragraph into the middle of \"Pristjek220\" at the end of the\n// administration-application paragraph.\n\nconst body = context.document.body;\n\n// --- 1. Small wording fixes -------------------------------------------------\nconst replacements = [\n  [\"handle. Derefter\", \"handle i. Derefter\"],\n  [\"Derudover kan man inde i den genereret indk\u00f8bsliste, skifte hvilken\", \"Der kan ogs\u00e5 inde i den genereret indk\u00f8bsliste, skiftes hvilken\"],\n  [\"hver enkelt produkt henne. Den genereret indk\u00f8bsliste viser ogs\u00e5, hvor meget bruger spare, i stedet for at k\u00f8be\", \"hver enkelt produkt i. Den genereret indk\u00f8bsliste viser samtidigt, hvor meget bruger spare, i forhold til at k\u00f8be\"],\n  [\"administratordel.\", \"administrator del.\"],\n  [\"administrere man en butik\", \"administrere man \u00e9n butik\"],\n  [\"kommer der et login s\u00e5 man skal\", \"kommer der et login hvor man skal\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${find}\", found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Move the `_GoBack` bookmark -----------------------------------------\n// Remove it from wherever it currently sits...\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ...and re-insert it inside \"Pristjek220.\" (between \"Pri\" and \"stjek220.\") at\n// the end of the administration-application paragraph.\nconst target = body.search(\"til Pri\", { matchCase: true });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length !== 1) {\n  throw new Error(\n    `Expected exactly 1 match for \"til Pri\", found ${target.items.length}`\n  );\n}\nconst caret = target.items[0].getRange(\"End\");\ncaret.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Rettet en smule i Resum\u00e9.docx\n# Applies the wording tweaks to the two \"Pristjek220\" summary paragraphs and\n# relocates the `_GoBack` bookmark from the end of the \"Gruppen har arbejdet...\"\n# paragraph into the middle of \"Pristjek220\" at the end of the\n# administration-application paragraph.\n\n$d = $word.ActiveDocument\n\n# wdReplaceOne = 1. Re-fetches $d.Content each call so earlier edits don't\n# invalidate the range, and fails loudly if a target string isn't found\n# (rather than silently leaving the document unchanged).\nfunction ReplaceOnce($findText, $replaceText) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $found) {\n        throw (\"Find.Execute did not find: \" + $findText)\n    }\n}\n\nReplaceOnce 'handle. Derefter' 'handle i. Derefter'\nReplaceOnce 'Derudover kan man inde i den genereret indk\u00f8bsliste, skifte hvilken' 'Der kan ogs\u00e5 inde i den genereret indk\u00f8bsliste, skiftes hvilken'\nReplaceOnce 'hver enkelt produkt henne. Den genereret indk\u00f8bsliste viser ogs\u00e5, hvor meget bruger spare, i stedet for at k\u00f8be' 'hver enkelt produkt i. Den genereret indk\u00f8bsliste viser samtidigt, hvor meget bruger spare, i forhold til at k\u00f8be'\nReplaceOnce 'administratordel.' 'administrator del.'\nReplaceOnce 'administrere man en butik' 'administrere man \u00e9n butik'\nReplaceOnce 'kommer der et login s\u00e5 man skal' 'kommer der et login hvor man skal'\n\n# --- Move the `_GoBack` bookmark --------------------------------------------\n# Drop it from wherever it currently sits...\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# ...and re-insert it inside \"Pristjek220.\" (between \"Pri\" and \"stjek220.\") at\n# the end of the administration-application paragraph.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"til Pri\", $false, $false, $false, $false, $false, $true, 1, $false)\nif (-not $found) {\n    throw \"Find.Execute did not find: til Pri\"\n}\n$rng.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n"}
